$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Kafka demo topic used throughout the cheatsheet
$ws.UsedRange.Replace("moon-landings", "lunar-landings") | Out-Null

# Add two new rows documenting how to list topics / consumer groups
$ws.Range("A11").Value = "Seznam topiců na brokeru"
$ws.Range("B11").Value = "kafka-topics --bootstrap-server=`$BOOT --list"
$ws.Range("A12").Value = "Seznam všech existujících consumer groups"
$ws.Range("B12").Value = "kafka-consumer-groups --bootstrap-server `$BOOT --list"

# Match styling/row height of the similarly-formatted rows above
$ws.Range("A11:A12").Font.WrapText = $true
$ws.Rows.Item(11).RowHeight = 17
$ws.Rows.Item(12).RowHeight = 17

$ws.Range("B6").Select() | Out-Null
